$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 67
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 1600
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44855
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("S3").Value = 3000
$ws.Range("T3").Value = 5

# Row 4
$ws.Range("D4").Value = 44875
$ws.Range("M4").Value = 50
